$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Ost/Nord rounded to nearest integer; Starttid (Z) / Sluttid (AB) cleared
$ws.Range("Q7").Value = 412577
$ws.Range("R7").Value = 6656304
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()

# Row 8 now carries what used to be row 13's species data
$ws.Range("A8").Value = 112083110
$ws.Range("B8").Value = 78107
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6453
$ws.Range("F8").Value = "Vedskivlav"
$ws.Range("G8").Value = "Hertelidea botryosa"
$ws.Range("H8").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q8").Value = 412206
$ws.Range("R8").Value = 6656051
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()

# Row 9: species data unchanged, only Ost/Nord rounded
$ws.Range("Q9").Value = 413016
$ws.Range("R9").Value = 6656415
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()

# Row 10 now carries what used to be row 11's species data
$ws.Range("A10").Value = 112083126
$ws.Range("B10").Value = 78536
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 229497
$ws.Range("F10").Value = "Korallblylav"
$ws.Range("G10").Value = "Parmeliella triptophylla"
$ws.Range("H10").Value = "(Ach.) Müll.Arg."
$ws.Range("Q10").Value = 413017
$ws.Range("R10").Value = 6656342
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()

# Row 11 now carries what used to be row 12's species data
$ws.Range("A11").Value = 112083127
$ws.Range("B11").Value = 77604
$ws.Range("E11").Value = 6450
$ws.Range("F11").Value = "Skuggblåslav"
$ws.Range("G11").Value = "Hypogymnia vittata"
$ws.Range("H11").Value = "(Ach.) Parrique"
$ws.Range("Q11").Value = 413052
$ws.Range("R11").Value = 6656343
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()

# Row 12 now carries what used to be row 8's species data
$ws.Range("A12").Value = 112083111
$ws.Range("B12").Value = 90666
$ws.Range("E12").Value = 4364
$ws.Range("F12").Value = "Dropptaggsvamp"
$ws.Range("G12").Value = "Hydnellum ferrugineum"
$ws.Range("H12").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q12").Value = 412205
$ws.Range("R12").Value = 6655989
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()

# Row 13 now carries what used to be row 14's species data
$ws.Range("A13").Value = 112083112
$ws.Range("B13").Value = 79444
$ws.Range("E13").Value = 1049
$ws.Range("F13").Value = "Kortskaftad ärgspik"
$ws.Range("G13").Value = "Microcalicium ahlneri"
$ws.Range("H13").Value = "Tibell"
$ws.Range("Q13").Value = 412284
$ws.Range("R13").Value = 6656072
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()

# Row 14 now carries what used to be row 10's species data
$ws.Range("A14").Value = 112083128
$ws.Range("B14").Value = 77186
$ws.Range("E14").Value = 353
$ws.Range("F14").Value = "Dvärgbägarlav"
$ws.Range("G14").Value = "Cladonia parasitica"
$ws.Range("H14").Value = "(Hoffm.) Hoffm."
$ws.Range("Q14").Value = 413190
$ws.Range("R14").Value = 6656475
$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()
